# Update viewer/ticket counts (column F) across sheets to reflect the
# latest scrape output, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" (Exhibitions) --
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 297
$ws1.Range("F5").Value = 1095
$ws1.Range("F6").Value = 3256
$ws1.Range("F10").Value = 712
$ws1.Range("F16").Value = 1510
$ws1.Range("F17").Value = 1510
$ws1.Range("F18").Value = 12
$ws1.Range("F20").Value = 25
$ws1.Range("F25").Value = 43385
$ws1.Range("F26").Value = 43385
$ws1.Range("F29").Value = 32239
$ws1.Range("F30").Value = 32239
$ws1.Range("F38").Value = 481
$ws1.Range("F39").Value = 1140
$ws1.Range("F40").Value = 5291
$ws1.Range("F45").Value = 310

# -- Sheet "演出" (Performances) --
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 14
$ws2.Range("F15").Value = 742
$ws2.Range("F22").Value = 0
$ws2.Range("F35").Value = 1132

# -- Sheet "本地生活" (Local life) --
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 520

# -- Sheet "全部类型" (All types) --
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 297
$ws4.Range("F4").Value = 520
$ws4.Range("F9").Value = 1095
$ws4.Range("F10").Value = 3256
$ws4.Range("F14").Value = 712
$ws4.Range("F23").Value = 1510
$ws4.Range("F24").Value = 1510
$ws4.Range("F27").Value = 25
$ws4.Range("F32").Value = 43385
$ws4.Range("F36").Value = 32239
$ws4.Range("F41").Value = 481
$ws4.Range("F42").Value = 1140
$ws4.Range("F43").Value = 5291
$ws4.Range("F49").Value = 310
